$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "Unknown" to "unknown" for the row-2 data cells D2:K2
$ws.Range("D2:K2").Value = "unknown"

$wb.Save()
